$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set A2 with the new label (no special formatting - default style)
$ws.Range("A2").Value = "param_pv1_area"

# Set B2 with the numeric value (no special formatting - default style)
$ws.Range("B2").Value = 2300
